# Dashboard Requirements.docx edit:
#  - Remove the "Total Delivered Items" and "Total Customers Analyzed" data-card
#    paragraphs from the Dashboard 2 (Customer) KPI list.
#  - Remove the "Total Delivered Items" and "Total Drivers Analyzed" data-card
#    paragraphs from the Dashboard 3 (Driver) KPI list.
#  - Add two more blank paragraphs to the existing run of blank paragraphs that
#    sits right before the "Dashboard 3: Drivers Loss Analysis" heading.

$d = $word.ActiveDocument

function Get-ParagraphIndexByText($searchText, $occurrenceIndex) {
    # Returns the 1-based Paragraphs index of the paragraph that contains the
    # Nth (occurrenceIndex) occurrence of searchText in the document body.
    $rng = $d.Content
    $n = 0
    while ($rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $n = $n + 1
        if ($n -eq $occurrenceIndex) {
            $foundStart = $rng.Start
            for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
                $p = $d.Paragraphs.Item($i)
                if ($foundStart -ge $p.Range.Start -and $foundStart -lt $p.Range.End) {
                    return $i
                }
            }
            return -1
        }
        $rng.Collapse(0)
    }
    return -1
}

function Remove-ParagraphByText($searchText, $occurrenceIndex) {
    # Deletes the whole paragraph (including its paragraph mark) that holds
    # the Nth occurrence of searchText.
    $idx = Get-ParagraphIndexByText $searchText $occurrenceIndex
    if ($idx -eq -1) {
        Write-Host "NOT FOUND: $searchText occurrence $occurrenceIndex"
        return
    }
    $p = $d.Paragraphs.Item($idx)
    $p.Range.Delete()
}

# --- Dashboard 2: Customer Loss Analysis data cards -----------------------
# "Total Delivered Items" appears 3 times overall (Dashboard 1, 2, 3); the
# Dashboard-1 one (1st) must stay, so remove the 2nd occurrence here.
Remove-ParagraphByText "Total Delivered Items" 2
Remove-ParagraphByText "Total Customers Analyzed" 1

# --- Dashboard 3: Drivers Loss Analysis data cards -------------------------
# After the removal above, the remaining "Total Delivered Items" text still
# occurs twice (Dashboard 1 and Dashboard 3); remove the 2nd (Dashboard 3) one.
Remove-ParagraphByText "Total Delivered Items" 2
Remove-ParagraphByText "Total Drivers Analyzed" 1

# --- Add two blank paragraphs before "Dashboard 3: ..." heading -----------
$idxHeading = Get-ParagraphIndexByText "Dashboard 3:" 1
$lastBlankPara = $d.Paragraphs.Item($idxHeading - 1)
$lastBlankPara.Range.InsertParagraphAfter()

$idxHeading = Get-ParagraphIndexByText "Dashboard 3:" 1
$lastBlankPara = $d.Paragraphs.Item($idxHeading - 1)
$lastBlankPara.Range.InsertParagraphAfter()
